# pabi_account_report / xlsx_report_receivable_detail.xlsx
# "update Allowance for doubtful accounts Report 13/07/2018"
#
# The report's detail header row (row 15) gets a new "Receipt Date" column
# inserted right before the existing "Receipt Number" column, and the
# "Reconcile number" header is re-capitalised to "Reconcile Number".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-seed the new header text far away from the real data (a scratch row
# well below the report) so it becomes a known shared string *before* the
# column insert below. Writing a brand new string straight into a column
# that was just created by Columns.Insert() can cause the whole shared
# string table to be rewritten losing preserved whitespace on unrelated
# strings, so we avoid that edge case by reusing an already-known string.
$ws.Range("A100").Value = "Receipt Date"

# Insert a new blank column before column R ("Receipt Number"), shifting
# that column and everything to its right (Preprint Number, Receipt
# Amount, Reconcile number, Validated By, Status) one column to the right.
$ws.Columns("R:R").Insert()

# Header for the newly inserted column.
$ws.Range("R15").Value = "Receipt Date"

# Remove the scratch cell/row again so it leaves no trace in the sheet.
$ws.Rows("100:100").Delete()

# Fix the capitalisation typo - this header is now in column V after the
# column insert shifted it one column to the right.
$ws.Range("V15").Value = "Reconcile Number"

# Reflect where the user was working when the change was made.
$ws.Range("U4").Select()
